$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.332.47"
$ws.Cells.Item(2, 5).Value = "  -3.21%  "

$ws.Cells.Item(3, 4).Value = "1.856.83"
$ws.Cells.Item(3, 5).Value = "  -3.93%  "

$ws.Cells.Item(4, 4).Value = "1.002"
$ws.Cells.Item(4, 5).Value = "  -0.24%  "

$ws.Cells.Item(5, 4).Value = "323.90"
$ws.Cells.Item(5, 5).Value = "  -2.65%  "

$ws.Cells.Item(6, 4).Value = "1.002"
$ws.Cells.Item(6, 5).Value = "  -0.21%  "

$ws.Cells.Item(7, 4).Value = "0.4529"
$ws.Cells.Item(7, 5).Value = "  -4.33%  "

$ws.Cells.Item(8, 4).Value = "0.3868"
$ws.Cells.Item(8, 5).Value = "  -4.78%  "

$ws.Cells.Item(9, 4).Value = "48.25"
$ws.Cells.Item(9, 5).Value = "  -8.99%  "

$ws.Cells.Item(10, 4).Value = "0.07930"
$ws.Cells.Item(10, 5).Value = "  -6.40%  "

$ws.Cells.Item(11, 5).Value = "  -3.63%  "

$ws.Cells.Item(12, 4).Value = "21.38"
$ws.Cells.Item(12, 5).Value = "  -4.22%  "

$ws.Cells.Item(13, 4).Value = "1.855.72"
$ws.Cells.Item(13, 5).Value = "  -5.65%  "

$ws.Cells.Item(14, 4).Value = "5.913"
$ws.Cells.Item(14, 5).Value = "  -3.51%  "

$ws.Cells.Item(15, 4).Value = "7.127"
$ws.Cells.Item(15, 5).Value = "  -5.49%  "

$ws.Cells.Item(16, 5).Value = "  -0.44%  "

$ws.Cells.Item(17, 4).Value = "85.89"
$ws.Cells.Item(17, 5).Value = "  -4.92%  "

$ws.Cells.Item(18, 5).Value = "  -3.83%  "

$ws.Cells.Item(19, 4).Value = "0.06548"

$ws.Cells.Item(20, 4).Value = "17.12"
$ws.Cells.Item(20, 5).Value = "  -6.21%  "

$ws.Cells.Item(21, 4).Value = "1.003"
$ws.Cells.Item(21, 5).Value = "  -0.19%  "

$ws.Cells.Item(22, 4).Value = "5.546"
$ws.Cells.Item(22, 5).Value = "  -4.22%  "

$ws.Cells.Item(23, 4).Value = "27.328.46"
$ws.Cells.Item(23, 5).Value = "  -3.34%  "

$ws.Cells.Item(24, 4).Value = "10.90"
$ws.Cells.Item(24, 5).Value = "  -4.92%  "

$ws.Cells.Item(25, 4).Value = "2.291"
$ws.Cells.Item(25, 5).Value = "  -0.09%  "

$ws.Cells.Item(26, 4).Value = "2.077.18"
$ws.Cells.Item(26, 5).Value = "  -5.31%  "

$ws.Cells.Item(27, 4).Value = "153.61"
$ws.Cells.Item(27, 5).Value = "  -0.54%  "

$ws.Cells.Item(28, 4).Value = "19.87"
$ws.Cells.Item(28, 5).Value = "  -1.44%  "

$ws.Cells.Item(29, 4).Value = "2.069"
$ws.Cells.Item(29, 5).Value = "  -4.66%  "

$ws.Cells.Item(30, 4).Value = "5.450"
$ws.Cells.Item(30, 5).Value = "  -5.80%  "

$ws.Cells.Item(31, 4).Value = "121.09"
$ws.Cells.Item(31, 5).Value = "  -2.19%  "

$ws.Cells.Item(32, 4).Value = "1.480"
$ws.Cells.Item(32, 5).Value = "  +1.51%  "

$ws.Cells.Item(33, 4).Value = "0.09292"
$ws.Cells.Item(33, 5).Value = "  -3.39%  "

$ws.Cells.Item(34, 4).Value = "0.9357"
$ws.Cells.Item(34, 5).Value = "  -5.04%  "

$ws.Cells.Item(35, 5).Value = "  -1.12%  "

$ws.Cells.Item(36, 4).Value = "5.269"
$ws.Cells.Item(36, 5).Value = "  -5.86%  "

$ws.Cells.Item(37, 4).Value = "1.227"
$ws.Cells.Item(37, 5).Value = "  -1.52%  "

$ws.Cells.Item(38, 5).Value = "  -4.25%  "

$ws.Cells.Item(39, 4).Value = "0.05986"
$ws.Cells.Item(39, 5).Value = "  -3.18%  "

$ws.Cells.Item(40, 4).Value = "8.130"
$ws.Cells.Item(40, 5).Value = "  -11.80%  "

$ws.Cells.Item(41, 4).Value = "1.001"
$ws.Cells.Item(41, 5).Value = "  -0.23%  "

$ws.Cells.Item(42, 4).Value = "0.5916"
$ws.Cells.Item(42, 5).Value = "  -4.59%  "

$ws.Cells.Item(43, 4).Value = "0.1888"
$ws.Cells.Item(43, 5).Value = "  -0.97%  "

$ws.Cells.Item(44, 4).Value = "10.14"
$ws.Cells.Item(44, 5).Value = "  -8.95%  "

$ws.Cells.Item(45, 4).Value = "1.273"
$ws.Cells.Item(45, 5).Value = "  -3.21%  "

$ws.Cells.Item(46, 4).Value = "0.5623"
$ws.Cells.Item(46, 5).Value = "  -4.84%  "

$ws.Cells.Item(47, 4).Value = "11.95"
$ws.Cells.Item(47, 5).Value = "  -7.57%  "

$ws.Cells.Item(48, 4).Value = "3.376"
$ws.Cells.Item(48, 5).Value = "  -3.01%  "

$ws.Cells.Item(49, 4).Value = "1.918"
$ws.Cells.Item(49, 5).Value = "  -6.50%  "

$ws.Cells.Item(50, 4).Value = "0.06752"
$ws.Cells.Item(50, 5).Value = "  -0.86%  "

$ws.Cells.Item(51, 4).Value = "108.75"
$ws.Cells.Item(51, 5).Value = "  -1.11%  "
